$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.988.30'
$ws.Range("E2").Value = '  +1.57%  '

# Row 3
$ws.Range("D3").Value = '3.258.47'
$ws.Range("E3").Value = '  +0.29%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.31'
$ws.Range("E5").Value = '  +1.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.93'
$ws.Range("E6").Value = '  +4.47%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  -1.10%  '

# Row 9
$ws.Range("E9").Value = '  +4.28%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").Value = '  -0.36%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.417'
$ws.Range("E11").Value = '  +2.07%  '

# Row 12
$ws.Range("D12").Value = '3.823.83'
$ws.Range("E12").Value = '  +0.41%  '

# Row 13
$ws.Range("E13").Value = '  +0.55%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.58'
$ws.Range("E14").Value = '  +2.19%  '

# Row 15
$ws.Range("D15").Value = '67.960.45'
$ws.Range("E15").Value = '  +1.56%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000173'
$ws.Range("E16").Value = '  +3.08%  '

# Row 17
$ws.Range("D17").Value = '3.261.94'
$ws.Range("E17").Value = '  +0.40%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.86'
$ws.Range("E18").Value = '  +1.00%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.61'
$ws.Range("E19").Value = '  +1.57%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '382.54'
$ws.Range("E20").Value = '  +3.43%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.68'
$ws.Range("E21").Value = '  +1.29%  '

# Row 22
$ws.Range("E22").Value = '  +0.00%  '

# Row 23
$ws.Range("E23").Value = '  +0.85%  '

# Row 24
$ws.Range("E24").Value = '  +1.12%  '

# Row 25
$ws.Range("E25").Value = '  +1.60%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.88'
$ws.Range("E26").Value = '  +1.10%  '

# Row 27
$ws.Range("E27").Value = '  +1.87%  '

# Row 28
$ws.Range("E28").Value = '  +0.00%  '

# Row 29
$ws.Range("E29").Value = '  +0.50%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.71'
$ws.Range("E30").Value = '  +1.35%  '

# Row 31
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.26'
$ws.Range("E31").Value = '  +7.29%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.96'
$ws.Range("E32").Value = '  +1.89%  '

# Row 34
$ws.Range("E34").Value = '  +3.27%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +3.29%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.10'
$ws.Range("E36").Value = '  -7.41%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.86'
$ws.Range("E37").Value = '  +0.83%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.837'
$ws.Range("E38").Value = '  -2.27%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.71'
$ws.Range("E39").Value = '  -0.33%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.71'
$ws.Range("E40").Value = '  +4.68%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.62'
$ws.Range("E41").Value = '  +7.45%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.61'
$ws.Range("E42").Value = '  +1.90%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.56'
$ws.Range("E43").Value = '  +3.67%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.37'
$ws.Range("E44").Value = '  +2.68%  '

# Row 45
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0690'
$ws.Range("E45").Value = '  +2.47%  '

# Row 46
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '345.64'
$ws.Range("E46").Value = '  +2.47%  '

# Row 47
$ws.Range("D47").Value = '2.645.66'
$ws.Range("E47").Value = '  -3.93%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0286'
$ws.Range("E48").Value = '  +2.85%  '

# Row 49
$ws.Range("E49").Value = '  -0.57%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("E50").Value = '  +1.96%  '

# Row 51
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.34'
$ws.Range("E51").Value = '  +3.08%  '
